# Apply the edit described by the diff:
#  - Remove the "mtn" data row (row 2) entirely, shifting the remaining
#    rows up by one (this also removes the now-unused shared strings for
#    the mtn url / mtn scraping pattern, and renumbers hyperlink rels).
#  - Update the interval (column B) for the "약업신문" (yakup) row from 6 to 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 2 (company = "mtn"), shifting rows 3-8 up to 2-7.
$ws.Rows.Item(2).Delete()

# After the shift, the row that used to be "약업신문" (row 5) is now row 4.
# Update its interval value (column B) from 6 to 9.
$ws.Cells.Item(4, 2).Value = 9
